$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gained a new "2020" data column (Q), extending the existing
# year series (D:P = 2007..2019) by one more year for every region row.

# 1) Copy the formatting (number format / font / border / alignment) of
#    the existing last data column (P) one column to the right (Q) so the
#    new column visually matches the rest of the table.
$ws.Range("P3:P14").Copy()
$ws.Range("Q3:Q14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 2) Fill in the new column's values: the 2020 header plus each region's
#    reported figure.
$ws.Range("Q4").Value = 2020
$ws.Range("Q5").Value = 38.6
$ws.Range("Q6").Value = 42.4
$ws.Range("Q7").Value = 53.2
$ws.Range("Q8").Value = 90.6
$ws.Range("Q9").Value = 52.6
$ws.Range("Q10").Value = 24.5
$ws.Range("Q11").Value = 69.1
$ws.Range("Q12").Value = 32.2
$ws.Range("Q13").Value = 19.1
$ws.Range("Q14").Value = 25.2

# 3) Leave the cursor where the author's session ended up.
$ws.Range("R27").Select()
